# Q3 Update - 2025
# Applies the data refresh to the UN-BDI (Burundi) "fromCSV" sheet:
#  - refresh the dataset-wide short-url token
#  - update a few existing rows' refugee/asylum-seeker figures
#  - correct the "Stateless" row (337) to actually be "South Sudan"
#  - insert a new "Stateless" row (338) with last quarter's figures
#  - push the old Uganda row down to 339 and bump its "items" counter

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) short-url column (B) changed for every data row: g0Xx6k -> IIeId4
# ------------------------------------------------------------------
$lastRow = $ws.Cells(1048576, 2).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 338 }
$ws.Range("B2:B$lastRow").Value = "IIeId4"

# ------------------------------------------------------------------
# 2) Row 330 - Burundi (IDPs within Burundi), 2024
# ------------------------------------------------------------------
$ws.Range("Q330").Value = "6877"
$ws.Range("T330").Value = "1281"

# ------------------------------------------------------------------
# 3) Row 332 - Dem. Rep. of the Congo -> Burundi, 2024
# ------------------------------------------------------------------
$ws.Range("N332").Value = "88199"
$ws.Range("O332").Value = "1946"
$ws.Range("P332").Value = "5"

# ------------------------------------------------------------------
# 4) Row 335 - Rwanda -> Burundi, 2024
# ------------------------------------------------------------------
$ws.Range("N335").Value = "823"
$ws.Range("O335").Value = "143"

# ------------------------------------------------------------------
# 5) Row 337 - was mislabeled "Stateless"; it is actually South Sudan
# ------------------------------------------------------------------
$ws.Range("F337").Value = "179"
$ws.Range("G337").Value = "South Sudan"
$ws.Range("H337").Value = "SSD"
$ws.Range("I337").Value = "SSD"
$ws.Range("N337").Value = "7"
$ws.Range("S337").Value = "0"

# ------------------------------------------------------------------
# 6) Insert a fresh row 338 for the real "Stateless" origin record
#    (pushes the old row 338 "Uganda" down to 339)
# ------------------------------------------------------------------
$ws.Rows(338).Insert()

$ws.Range("A338").Value = "1"
$ws.Range("B338").Value = "IIeId4"
$ws.Range("C338").Value = "1"
$ws.Range("D338").Value = "337"
$ws.Range("E338").Value = "2024"
$ws.Range("F338").Value = "216"
$ws.Range("G338").Value = "Stateless"
$ws.Range("H338").Value = "STA"
$ws.Range("I338").Value = "XXA"
$ws.Range("J338").Value = "16"
$ws.Range("K338").Value = "Burundi"
$ws.Range("L338").Value = "BDI"
$ws.Range("M338").Value = "BDI"
$ws.Range("N338").Value = "0"
$ws.Range("O338").Value = "0"
$ws.Range("P338").Value = "0"
$ws.Range("Q338").Value = "0"
$ws.Range("R338").Value = "0"
$ws.Range("S338").Value = "791"
$ws.Range("T338").Value = "0"
$ws.Range("U338").Value = "-"
$ws.Range("V338").Value = "0"

# ------------------------------------------------------------------
# 7) Row 339 (the shifted Uganda row) gets a bumped "items" index
# ------------------------------------------------------------------
$ws.Range("D339").Value = "338"
